$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ChaylaGreyJourney")

$ws.Range("B4").Value = "cum with me right now papi... watch every fucking second 😏"
$ws.Range("B6").Value = "I'm cumming papi... fuck, watch me let go all over for you"
$ws.Range("B7").Value = "I'm right there... don't you dare cum before you watch me finish first 😏"
$ws.Range("B8").Value = "oh my fucking god 🔥"
$ws.Range("B9").Value = "you're about to see what happens when I completely lose control... this is all you 😏"
$ws.Range("B11").Value = "I'm about to lose it papi... you need to watch what you did to me"
$ws.Range("B12").Value = "I want to ride you so bad while you grab my hips and don't let go... I'm losing my mind 😏"
$ws.Range("B13").Value = "I'm playing with my pussy right now imagining you deep inside me... I need to feel every inch"
$ws.Range("B14").Value = "FUCK 🔥"
$ws.Range("B15").Value = "see what you're doing to me papi? I can't stop and I don't want to 😏"
$ws.Range("B17").Value = "tell me exactly how you want me... I'll do whatever you say right now"
$ws.Range("B18").Value = "I need your hands all over my body so bad it almost hurts papi... feel how wet you're making me 🔥"
$ws.Range("B19").Value = "I'm dripping wet right now thinking about what I want to do to you... god I need it"
$ws.Range("B20").Value = "mm that was just the warmup 😏"
$ws.Range("B21").Value = "look at what you started... hope you can handle this 😏"
$ws.Range("B23").Value = "I'm already touching myself and it's your fault papi... hope you can handle what comes next 😏"
$ws.Range("B24").Value = "talking to you is making me so turned on right now... I can feel it building and I'm done holding back"
$ws.Range("B25").Value = "you like what you see? because now I'm really in the mood to show you more 🔥"
